$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "color"
$ws.Range("D2").Value = "White"
$ws.Range("D3").Value = "Black"
$ws.Range("D4").Value = "White"
$ws.Range("D5").Value = "Black"
$ws.Range("D6").Value = "Green"
$ws.Range("D7").Value = "Red"
$ws.Range("D8").Value = "Yellow"
$ws.Range("D9").Value = "Purple"
$ws.Range("D10").Value = "Purple"
$ws.Range("D11").Value = "Green"

$ws.Range("D11").Select()
